$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) were re-shuffled: each row now carries the values
# that (in the original workbook) belonged to a different row, for
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg).
#
# Target values after the edit:
$rows = @(
    @{ Row = 2; D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 },
    @{ Row = 3; D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 },
    @{ Row = 4; D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 },
    @{ Row = 5; D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 },
    @{ Row = 6; D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("P$row").Value = $r.P
}
